$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: update the "Price" (column D) and "Volume(1h)"
# (column E) figures for each coin row, keeping everything else untouched.
#
# A few new Price values (e.g. "1.000") would otherwise be auto-detected
# by Excel as numbers and lose their trailing zeros, so those specific
# cells are written with a leading apostrophe to force plain text, just
# like typing them directly into the grid.

$ws.Range("D2").Value = "27.760.78"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.743.61"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "333.03"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "0.3879"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").Value = "0.3363"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "45.33"
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("D10").Value = "1.097"
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").Value = "0.07127"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "21.75"
$ws.Range("E13").Value = "  -6.22%  "
$ws.Range("D14").Value = "6.062"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Value = "1.742.08"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "6.929"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "0.00001044"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "0.06594"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "78.72"
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  -5.15%  "
$ws.Range("D22").Value = "6.153"
$ws.Range("E22").Value = "  -4.55%  "
$ws.Range("D23").Value = "27.749.64"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("D25").Value = "2.389"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "153.38"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "19.68"
$ws.Range("E27").Value = "  -6.25%  "
$ws.Range("D28").Value = "2.268"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("D29").Value = "1.940.37"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "1.272"
$ws.Range("E30").Value = "  -11.52%  "
$ws.Range("D31").Value = "126.94"
$ws.Range("E31").Value = "  -6.13%  "
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("D33").Value = "5.724"
$ws.Range("E33").Value = "  -7.51%  "
$ws.Range("D34").Value = "0.08684"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "11.92"
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("D36").Value = "1.505"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "5.069"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D38").Value = "0.02248"
$ws.Range("E38").Value = "  -7.76%  "
$ws.Range("D39").Value = "0.06044"
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("D40").Value = "0.6379"
$ws.Range("E40").Value = "  -7.51%  "
$ws.Range("D41").Value = "0.2073"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("D42").Value = "1.189"
$ws.Range("E42").Value = "  -4.53%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "7.803"
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "3.804"
$ws.Range("D47").Value = "0.5882"
$ws.Range("E47").Value = "  -7.28%  "
$ws.Range("D48").Value = "125.19"
$ws.Range("E48").Value = "  -5.46%  "
$ws.Range("D49").Value = "'1.960"
$ws.Range("E49").Value = "  -6.88%  "
$ws.Range("D50").Value = "0.06919"
$ws.Range("E50").Value = "  -7.56%  "
$ws.Range("E51").Value = "  -5.10%  "
